$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Extend the existing year columns (B..J = 2014..2022) with a new
# column K for year 2023, mirroring the formatting of the preceding
# column J for each of the three data rows.
$ws.Range("J3:J6").Copy()
$ws.Range("K3:K6").PasteSpecial(-4122)

$ws.Range("K3").Value = 2023
$ws.Range("K4").Value = 410.6
$ws.Range("K5").Value = 373.2
$ws.Range("K6").Value = 425.3
